$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: expand the table to its final size BEFORE inserting rows so the
#     structured-reference formulas keep the Table1[[#This Row],[...]] form
#     instead of collapsing to [@...] once rows temporarily sit outside ref.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K831"))

# --- Step 2: insert 3 new blank rows right after row 664 (i.e. rows 665:667),
#     pushing every following row down by 3 (828 -> 831 etc.)
$ws.Rows("665:667").Insert()

# --- Step 3: the plain row insert drops borders on the new rows; restore the
#     exact formatting by copying it from an untouched template row (670).
$ws.Range("A670:K670").Copy()
$ws.Range("A665:K667").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 4: the calculated "EARNED " column formula was wiped by the
#     formats-only paste on the 3 new rows; put it back.
$formula667 = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G665").Formula = $formula667
$ws.Range("G666").Formula = $formula667
$ws.Range("G667").Formula = $formula667

# --- Step 5: REMARKS (K) column on rows 665-668 uses the date-formatted
#     style (same as the existing K664) rather than the plain style the
#     template row carried.
$ws.Range("K664").Copy()
$ws.Range("K665:K668").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 6: populate the new leave-card entries.
# Row 664 (existing row): SL(1-0-0) now also shows an EARNED amount of 1.25.
$ws.Range("C664").Value = 1.25

# Row 665 (new): SL(2-0-0) taken, 2 day(s) charged to vacation absence col,
# with remarks referencing the leave dates.
$ws.Range("B665").Value = "SL(2-0-0)"
$ws.Range("H665").Value = 2
$ws.Range("K665").Value = "9/29 , 10/1/2023"

# Row 666 (new): SL(1-0-0) earned 1.25, with 1 day vacation absence and a
# remarks date.
$ws.Range("B666").Value = "SL(1-0-0)"
$ws.Range("C666").Value = 1.25
$ws.Range("H666").Value = 1
$ws.Range("K666").Value = 45209

# Row 667 (new): VL(2-0-0), 2 days sick-leave absence, remarks date string.
$ws.Range("B667").Value = "VL(2-0-0)"
$ws.Range("D667").Value = 2
$ws.Range("K667").Value = "10/23,24/2023"

# Row 668 (new): SL(1-0-0), 1 day vacation absence, remarks date.
$ws.Range("B668").Value = "SL(1-0-0)"
$ws.Range("H668").Value = 1
$ws.Range("K668").Value = 45224

# Row 669 (new): VL(3-0-0), 3 days sick-leave absence, remarks date string.
$ws.Range("B669").Value = "VL(3-0-0)"
$ws.Range("D669").Value = 3
$ws.Range("K669").Value = "12/27,28,29/2023"
